$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.661.22'
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").Value = '1.859.23'
$ws.Range("E3").Value = '  +0.39%  '

$c = $ws.Range("D4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.020'
$c.Style = $s
$ws.Range("E4").Value = '  -1.13%  '

$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '320.34'
$c.Style = $s
$ws.Range("E5").Value = '  -0.38%  '

$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.018'
$c.Style = $s
$ws.Range("E6").Value = '  -0.88%  '

$c = $ws.Range("D7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.4364'
$c.Style = $s
$ws.Range("E7").Value = '  -0.50%  '

$c = $ws.Range("D8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.3793'
$c.Style = $s
$ws.Range("E8").Value = '  +0.49%  '

$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.07421'
$c.Style = $s
$ws.Range("E9").Value = '  +0.14%  '

$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.8842'
$c.Style = $s
$ws.Range("E10").Value = '  +0.92%  '

$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '21.59'
$c.Style = $s
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").Value = '1.880.60'
$ws.Range("E12").Value = '  +1.62%  '

$c = $ws.Range("D13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.740'
$c.Style = $s
$ws.Range("E13").Value = '  +0.61%  '

$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.485'
$c.Style = $s
$ws.Range("E14").Value = '  -0.80%  '

$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.07094'
$c.Style = $s
$ws.Range("E15").Value = '  -1.60%  '

$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '86.97'
$c.Style = $s
$ws.Range("E16").Value = '  +4.91%  '

$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.024'
$c.Style = $s
$ws.Range("E17").Value = '  -1.03%  '

$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.000009060'
$c.Style = $s
$ws.Range("E18").Value = '  +0.39%  '

$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.019'
$c.Style = $s
$ws.Range("E19").Value = '  -0.92%  '

$c = $ws.Range("D20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '15.45'
$c.Style = $s
$ws.Range("E20").Value = '  +0.24%  '

$ws.Range("D21").Value = '27.722.78'
$ws.Range("E21").Value = '  +0.55%  '

$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.282'
$c.Style = $s
$ws.Range("E22").Value = '  +0.49%  '

$ws.Range("E23").Value = '  -1.60%  '

$ws.Range("D24").Value = '2.087.87'
$ws.Range("E24").Value = '  +0.76%  '

$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.035'
$c.Style = $s
$ws.Range("E25").Value = '  +6.32%  '

$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '157.47'
$c.Style = $s
$ws.Range("E26").Value = '  -0.22%  '

$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '18.71'
$c.Style = $s
$ws.Range("E27").Value = '  -0.34%  '

$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.371'
$c.Style = $s
$ws.Range("E28").Value = '  +2.06%  '

$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.997'
$c.Style = $s
$ws.Range("E29").Value = '  +1.23%  '

$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '120.42'
$c.Style = $s
$ws.Range("E30").Value = '  +2.91%  '

$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.09058'
$c.Style = $s
$ws.Range("E31").Value = '  +0.06%  '

$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.218'
$c.Style = $s
$ws.Range("E32").Value = '  +1.68%  '

$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.7688'
$c.Style = $s
$ws.Range("E33").Value = '  +0.83%  '

$ws.Range("E34").Value = '  +5.55%  '

$c = $ws.Range("D35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.558'
$c.Style = $s
$ws.Range("E35").Value = '  +0.98%  '

$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.019'
$c.Style = $s
$ws.Range("E36").Value = '  -0.88%  '

$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.144'
$c.Style = $s
$ws.Range("E37").Value = '  -0.39%  '

$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.01978'
$c.Style = $s
$ws.Range("E38").Value = '  +0.34%  '

$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.05295'
$c.Style = $s
$ws.Range("E39").Value = '  -0.11%  '

$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.875'
$c.Style = $s
$ws.Range("E40").Value = '  +2.66%  '

$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.5195'
$c.Style = $s
$ws.Range("E41").Value = '  +0.78%  '

$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.947'
$c.Style = $s
$ws.Range("E42").Value = '  +3.11%  '

$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.1682'
$c.Style = $s
$ws.Range("E43").Value = '  +0.40%  '

$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '8.696'
$c.Style = $s
$ws.Range("E44").Value = '  +2.53%  '

$ws.Range("E45").Value = '  +2.60%  '

$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '110.09'
$c.Style = $s
$ws.Range("E46").Value = '  +1.23%  '

$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.713'
$c.Style = $s
$ws.Range("E47").Value = '  +0.21%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.4727'
$c.Style = $s
$ws.Range("E48").Value = '  +1.68%  '

$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.019'
$c.Style = $s
$ws.Range("E49").Value = '  -1.14%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.06497'
$c.Style = $s
$ws.Range("E50").Value = '  +1.67%  '

$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.861'
$c.Style = $s
$ws.Range("E51").Value = '  +0.20%  '

